# "more reports from the past few days"
#
# 1. Append six more days (rows 98-103) of train-run data to the "Data"
#    sheet, matching the new daily completion / trip-length numbers.
# 2. Add the new blank "DE.1.0.7.0 comparisons" worksheet.
# 3. Leave the selection on the last pair of cells that were touched
#    (H102:H103), matching where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Columns: B=Scheduled, C=Completed, D=Skipped, E=Incomplete, F=Total Completed,
#          G=Total Completed %, H/I/J = run-time stats (avg/min/max minutes)
$newRows = @(
    @{ Row = 98;  B = 143; C = 133; D = 0; E = 10; F = 133; G = 0.93006993006993011;  H = 43.965037593896383; I = 35.333333326270804; J = 64.599999992642552 },
    @{ Row = 99;  B = 138; C = 132; D = 0; E = 6;  F = 132; G = 0.95652173913043481;  H = 46.900252525837544; I = 34.966666668187827; J = 107.68333333893679 },
    @{ Row = 100; B = 140; C = 137; D = 0; E = 3;  F = 137; G = 0.97857142857142854;  H = 42.51934306543103;  I = 34.800000004470348; J = 55.33333332859911  },
    @{ Row = 101; B = 140; C = 134; D = 0; E = 6;  F = 134; G = 0.95714285714285718;  H = 43.084999999655075; I = 35.016666669398546; J = 79.049999996786937 },
    @{ Row = 102; B = 145; C = 143; D = 0; E = 2;  F = 143; G = 0.98620689655172411;  H = 43.229836829566771; I = 35.416666668606922; J = 62.133333336096257 },
    @{ Row = 103; B = 138; C = 131; D = 0; E = 7;  F = 131; G = 0.94927536231884058;  H = 44.978371500455282; I = 36.083333323476836; J = 64.74999999627471  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
}

# New (still empty) worksheet for the next comparison report, placed at
# the end of the tab strip.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "DE.1.0.7.0 comparisons"

# Leave the Data sheet active with the last-edited cells selected, as in
# the author's saved view.
$ws.Activate()
$ws.Range("H102:H103").Select()
